$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply weekly price-sheet update: re-sorted rows by date, one row price revision,
# and a new trailing row (12) appended with the previously-missing entry.

# Row 2
$ws.Range("D2").Value = 44462
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N2").Value = 2900
$ws.Range("O2").Value = 3000
$ws.Range("P2").Value = 2950
$ws.Range("Q2").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S2").Value = 2950
$ws.Range("T2").Value = 1

# Row 3
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 2600
$ws.Range("O3").Value = 2600
$ws.Range("P3").Value = 2600
$ws.Range("S3").Value = 2600

# Row 4
$ws.Range("D4").Value = 44160
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = "$/bandeja 8 kilos"
$ws.Range("S4").Value = 2188
$ws.Range("T4").Value = 8

# Row 5
$ws.Range("D5").Value = 44160
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/bandeja 8 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 1875
$ws.Range("T5").Value = 8

# Row 6
$ws.Range("D6").Value = 44446
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N6").Value = 3200
$ws.Range("O6").Value = 3300
$ws.Range("P6").Value = 3250
$ws.Range("Q6").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R6").Value = "Provincia del Elquí"
$ws.Range("S6").Value = 3250
$ws.Range("T6").Value = 1

# Row 7
$ws.Range("D7").Value = 44454
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N7").Value = 30000
$ws.Range("O7").Value = 31000
$ws.Range("P7").Value = 30500
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("S7").Value = 3050
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 2700
$ws.Range("O8").Value = 2800
$ws.Range("P8").Value = 2750
$ws.Range("S8").Value = 2750

# Row 9
$ws.Range("D9").Value = 44467
$ws.Range("D9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 2500
$ws.Range("O9").Value = 2500
$ws.Range("P9").Value = 2500
$ws.Range("Q9").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S9").Value = 2500
$ws.Range("T9").Value = 1

# Row 10
$ws.Range("D10").Value = 44469
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 28000
$ws.Range("O10").Value = 29000
$ws.Range("P10").Value = 28500
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("S10").Value = 2850
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("D11").Value = 44475
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 12
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Vega Monumental Concepción"
$ws.Range("C12").Value = "Bíobío"
$ws.Range("D12").Value = 44461
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = "Otros"
$ws.Range("I12").Value = 100107002
$ws.Range("J12").Value = "Chirimoya"
$ws.Range("K12").Value = "Cultivar IV Región"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 29000
$ws.Range("O12").Value = 30000
$ws.Range("P12").Value = 29500
$ws.Range("Q12").Value = "$/bandeja 10 kilos"
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 2950
$ws.Range("T12").Value = 10
